$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "exo shirt",
    "exo shirts",
    "exo short",
    "exo suit",
    "exo suit real",
    "exo winter",
    "exos 38",
    "exos back brace",
    "exos form 2 back brace",
    "exos form back brace",
    "extra clothing",
    "extra large knee support",
    "extra long black leggings for tall women",
    "extra long compression shorts",
    "extra long compression tank tops for women",
    "extra long legging",
    "extra long running pants",
    "extra long sports pants",
    "extra long womens leggings",
    "extra support",
    "eye compression",
    "eye tight",
    "f34 cwx ss",
    "f40 cwx",
    "fabletics 3 pack of leggings",
    "fabletics capris for women",
    "fabletics leggings capri",
    "fabletics leggings women",
    "fabric frisbee",
    "fabric leggings for women",
    "fabric skin tape",
    "fabric to skin glue",
    "fabric to skin tape",
    "fabulous leggings",
    "fabulous workout clothes",
    "fall cycling pants",
    "fall recovery",
    "fashion capri leggings",
    "fashion capri leggings for women",
    "fashion pants",
    "fashion short tights for women",
    "fashion tights",
    "fashion to figure pants",
    "fast drying pants women",
    "fast exercise",
    "fast line",
    "fatigue compression pants",
    "fatigue leggings",
    "fatigue leggings for women",
    "fatigue pants for women",
    "fatigue runner",
    "fatigue shorts for women",
    "fatigue tights",
    "fatigued to fantastic",
    "fc16t9 cw",
    "features compression",
    "feeling chart for therapy",
    "felt balls pastel",
    "felt lined pants men",
    "felt pastel",
    "felt soccer ball",
    "felt soccer ball size 4",
    "felt wick",
    "felt winter",
    "female bike pants",
    "female capri pants",
    "female compression",
    "female compression leggings",
    "female compression pants",
    "female compression shirt",
    "female compression shorts",
    "female joint",
    "female pelvis model",
    "female running gear",
    "female running shorts",
    "female sports tights",
    "female support",
    "female tommies",
    "female training",
    "female wears",
    "female workout leggings",
    "fibromyalgia back brace",
    "fibromyalgia clothing",
    "fibromyalgia device",
    "fibromyalgia exercise",
    "fibromyalgia knee brace",
    "fibromyalgia relief ball",
    "fibromyalgia shirts",
    "fibromyalgia support",
    "ficus fig",
    "field and sport pro",
    "field leggings",
    "field runner",
    "field sport",
    "field target",
    "fig clothing women",
    "fila pants women",
    "final expense",
    "finals week gifts",
    "firm compression waist",
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
